# Add two new services, two new quotes, and two new quote_items rows
# mirroring the structure/style of existing rows in each sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# services sheet: append rows 19 and 20
# ---------------------------------------------------------------
$wsServices = $wb.Worksheets.Item("services")

$wsServices.Cells.Item(19, 1).Value = 18
$wsServices.Cells.Item(19, 2).Value = "eletrodomesticos"
$wsServices.Cells.Item(19, 4).Value = 47.4
$wsServices.Cells.Item(19, 7).Value = "unidade"
$wsServices.Cells.Item(19, 13).Value = "2025-09-25T13:53:26.215170"
$wsServices.Cells.Item(19, 14).Value = "2025-09-25T13:53:26.215170"

$wsServices.Cells.Item(20, 1).Value = 19
$wsServices.Cells.Item(20, 2).Value = "dedetizacao"
$wsServices.Cells.Item(20, 4).Value = 129.94
$wsServices.Cells.Item(20, 7).Value = "unidade"
$wsServices.Cells.Item(20, 13).Value = "2025-09-25T13:58:09.352426"
$wsServices.Cells.Item(20, 14).Value = "2025-09-25T13:58:09.352426"

# ---------------------------------------------------------------
# quotes sheet: append rows 13 and 14
# ---------------------------------------------------------------
$wsQuotes = $wb.Worksheets.Item("quotes")

$wsQuotes.Cells.Item(13, 1).Value = 13
$wsQuotes.Cells.Item(13, 2).Value = "ORC202509011"
$wsQuotes.Cells.Item(13, 3).Value = 1
$wsQuotes.Cells.Item(13, 4).Value = "Orçamento - eletrodomesticos"
$wsQuotes.Cells.Item(13, 5).Value = "instalacoes"
$wsQuotes.Cells.Item(13, 8).Value = "pendente"
$wsQuotes.Cells.Item(13, 13).Value = 47.4
$wsQuotes.Cells.Item(13, 18).Value = "2025-09-25T13:53:26.691730"
$wsQuotes.Cells.Item(13, 19).Value = "2025-09-25T13:53:26.691730"

$wsQuotes.Cells.Item(14, 1).Value = 14
$wsQuotes.Cells.Item(14, 2).Value = "ORC202509012"
$wsQuotes.Cells.Item(14, 3).Value = 1
$wsQuotes.Cells.Item(14, 4).Value = "Orçamento - dedetizacao"
$wsQuotes.Cells.Item(14, 5).Value = "limpeza"
$wsQuotes.Cells.Item(14, 8).Value = "pendente"
$wsQuotes.Cells.Item(14, 13).Value = 129.94
$wsQuotes.Cells.Item(14, 18).Value = "2025-09-25T13:58:10.103277"
$wsQuotes.Cells.Item(14, 19).Value = "2025-09-25T13:58:10.103277"

# ---------------------------------------------------------------
# quote_items sheet: append rows 14 and 15
# ---------------------------------------------------------------
$wsItems = $wb.Worksheets.Item("quote_items")

$wsItems.Cells.Item(14, 1).Value = 15
$wsItems.Cells.Item(14, 2).Value = 13
$wsItems.Cells.Item(14, 3).Value = 18
$wsItems.Cells.Item(14, 4).Value = 1
$wsItems.Cells.Item(14, 5).Value = 47.4
$wsItems.Cells.Item(14, 6).Value = ""
$wsItems.Cells.Item(14, 7).Value = 47.4
$wsItems.Cells.Item(14, 8).Value = "eletrodomesticos"
$wsItems.Cells.Item(14, 9).Value = ""
$wsItems.Cells.Item(14, 10).Value = "unidade"
$wsItems.Cells.Item(14, 11).Value = ""
$wsItems.Cells.Item(14, 12).Value = ""
$wsItems.Cells.Item(14, 13).Value = ""
$wsItems.Cells.Item(14, 14).Value = ""
$wsItems.Cells.Item(14, 15).Value = "2025-09-25T13:53:26.691730"

$wsItems.Cells.Item(15, 1).Value = 16
$wsItems.Cells.Item(15, 2).Value = 14
$wsItems.Cells.Item(15, 3).Value = 19
$wsItems.Cells.Item(15, 4).Value = 1
$wsItems.Cells.Item(15, 5).Value = 129.94
$wsItems.Cells.Item(15, 6).Value = ""
$wsItems.Cells.Item(15, 7).Value = 129.94
$wsItems.Cells.Item(15, 8).Value = "dedetizacao"
$wsItems.Cells.Item(15, 9).Value = ""
$wsItems.Cells.Item(15, 10).Value = "unidade"
$wsItems.Cells.Item(15, 11).Value = ""
$wsItems.Cells.Item(15, 12).Value = ""
$wsItems.Cells.Item(15, 13).Value = ""
$wsItems.Cells.Item(15, 14).Value = ""
$wsItems.Cells.Item(15, 15).Value = "2025-09-25T13:58:10.103277"
